$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 339; existing rows 339-429 shift down to 340-430.
$ws.Rows(339).Insert()

$ws.Cells.Item(339, 1).Value = 5
$ws.Cells.Item(339, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(339, 3).Value = "Maule"
$ws.Cells.Item(339, 4).Value = 45135
$ws.Cells.Item(339, 5).Value = 7
$ws.Cells.Item(339, 6).Value = "Fruta"
$ws.Cells.Item(339, 7).Value = 100108
$ws.Cells.Item(339, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(339, 9).Value = 100108005
$ws.Cells.Item(339, 10).Value = "Piña"
$ws.Cells.Item(339, 11).Value = "Sin especificar"
$ws.Cells.Item(339, 12).Value = "Segunda"
$ws.Cells.Item(339, 13).Value = 250
$ws.Cells.Item(339, 14).Value = 18000
$ws.Cells.Item(339, 15).Value = 18000
$ws.Cells.Item(339, 16).Value = 18000
$ws.Cells.Item(339, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(339, 18).Value = "Ecuador"
$ws.Cells.Item(339, 19).Value = 1286
$ws.Cells.Item(339, 20).Value = 14
